# planning_projet.xlsx edit
# Fill in the "Réalisée" (G) column for the sub-tasks of rows 7-26:
# most become "Oui", the freshly-introduced ones for row 7/8/9/11 become "Skipped".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("G7").Value = "Skipped"
$ws.Range("G8").Value = "Skipped"
$ws.Range("G9").Value = "Skipped"
$ws.Range("G10").Value = "Oui"
$ws.Range("G11").Value = "Skipped"

# Row 12 lost its explicit-black font in the real edit (cell style moved from
# the "fontId 3 / black" look to the plain theme-coloured centred look that
# D46/D47/D48/D50 already use on this sheet) - copy that existing format over
# before writing the value so we reuse the existing style record.
$ws.Range("D46").Copy() | Out-Null
$ws.Range("G12").PasteSpecial(-4122) | Out-Null
$ws.Range("G12").Value = "Oui"

$ws.Range("G13").Value = "Oui"
$ws.Range("G14").Value = "Oui"
$ws.Range("G15").Value = "Oui"
$ws.Range("G16").Value = "Oui"
$ws.Range("G17").Value = "Oui"
$ws.Range("G18").Value = "Oui"
$ws.Range("G19").Value = "Oui"
$ws.Range("G20").Value = "Oui"
$ws.Range("G21").Value = "Oui"
$ws.Range("G22").Value = "Oui"
$ws.Range("G24").Value = "Oui"
$ws.Range("G26").Value = "Oui"

# Mirror the author's final on-screen selection.
$ws.Range("G23").Select() | Out-Null
